$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was inserted as row 96, pushing the previous rows 96-131
# down to 97-132 (dimension grows from A1:R131 to A1:R132).
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new price record.
$ws.Range("A96").Value = 1
$ws.Range("B96").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C96").Value = "Arica y Parinacota"
$ws.Range("D96").Value = 44809
$ws.Range("E96").Value = 15
$ws.Range("F96").Value = 100112042
$ws.Range("G96").Value = "Locoto"
$ws.Range("H96").Value = "Sin especificar"
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 140
$ws.Range("K96").Value = 20000
$ws.Range("L96").Value = 21000
$ws.Range("M96").Value = 20500
$ws.Range("N96").Value = "`$/caja 20 kilos"
$ws.Range("O96").Value = "Región de Arica y Parinacota"
$ws.Range("P96").Value = 1025
$ws.Range("Q96").Value = 20
$ws.Range("R96").Value = "Hortaliza"
